$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.127.01'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.734.52'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '616.88'
$ws.Range("E5").Value = '  +6.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.84'
$ws.Range("E6").Value = '  +3.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.733.43'
$ws.Range("E7").Value = '  -0.72%  '
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.722'
$ws.Range("E10").Value = '  -1.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.162'
$ws.Range("E11").Value = '  -5.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '56.80'
$ws.Range("E12").Value = '  +5.97%  '
$ws.Range("E13").Value = '  -4.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.69'
$ws.Range("E14").Value = '  -2.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.323.49'
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.731.18'
$ws.Range("E16").Value = '  -2.46%  '
$ws.Range("E17").Value = '  -1.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.10'
$ws.Range("E18").Value = '  -1.75%  '
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("E20").Value = '  -2.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '68.975.97'
$ws.Range("E21").Value = '  -0.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '414.13'
$ws.Range("E22").Value = '  -1.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.69'
$ws.Range("E23").Value = '  +1.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.76'
$ws.Range("E24").Value = '  -1.10%  '
$ws.Range("E25").Value = '  -2.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.83'
$ws.Range("E26").Value = '  -3.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.95'
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("E28").Value = '  +2.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.81'
$ws.Range("E29").Value = '  -0.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.68'
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.27'
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.36'
$ws.Range("E32").Value = '  -12.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.74'
$ws.Range("E33").Value = '  -2.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.123'
$ws.Range("E34").Value = '  +1.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '620.94'
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '44.53'
$ws.Range("E36").Value = '  -2.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '66.06'
$ws.Range("E37").Value = '  -1.45%  '
$ws.Range("E38").Value = '  -7.82%  '
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("E40").Value = '  -0.98%  '
$ws.Range("E41").Value = '  -0.65%  '
$ws.Range("E42").Value = '  +1.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.07'
$ws.Range("E43").Value = '  -2.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0445'
$ws.Range("E44").Value = '  -1.07%  '
$ws.Range("E45").Value = '  +0.50%  '
$ws.Range("E46").Value = '  +1.04%  '
$ws.Range("E47").Value = '  -4.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.839.86'
$ws.Range("E48").Value = '  +1.89%  '
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.72'
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.68'
$ws.Range("E50").Value = '  -16.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.12'
$ws.Range("E51").Value = '  -3.83%  '
